$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 10; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Value = "cost_variable_om"
}

$ws.Range("C10:C39").Select()
